# Applies the "add bad gateway on server error" commit's document edits.
# The underlying change swaps a handful of placeholder/test values inside
# the student-output table report. Cells are addressed by
# Table/Row/Cell index so that the many duplicate "No"/"Yes" strings in
# the document are targeted unambiguously.

$d = $word.ActiveDocument

# --- Table 1: Company overview -------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.Rows.Item(1).Cells.Item(2).Range.Text = "Ashutosh Tripathi"   # Name of the Company
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "google.com"          # Website

# --- Table 2: Intern profile ----------------------------------------------
$tbl = $d.Tables.Item(2)
$tbl.Rows.Item(2).Cells.Item(2).Range.Text = "sde 21"               # Job Designation
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "sde 21"               # Job Description
$tbl.Rows.Item(5).Cells.Item(2).Range.Text = "sde 21"               # Place of Posting

# --- Table 3: Stipend / PPO / CTC -----------------------------------------
$tbl = $d.Tables.Item(3)
$tbl.Rows.Item(1).Cells.Item(2).Range.Text = "sde 21"               # Stipend per month
$tbl.Rows.Item(2).Cells.Item(2).Range.Text = "No"                   # PPO provision (Yes -> No)
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "sde 21"               # CTC for PPO selects

# --- Table 4: Courses checkboxes (SELECT ALL + every course: No -> Yes) --
$tbl = $d.Tables.Item(4)
$rowCount = $tbl.Rows.Count
for ($i = 2; $i -le $rowCount; $i++) {
    $tbl.Rows.Item($i).Cells.Item(2).Range.Text = "Yes"
}

# --- Table 10: Interview process summary ----------------------------------
$tbl = $d.Tables.Item(10)
$tbl.Rows.Item(9).Cells.Item(2).Range.Text = "ashutosh"             # Total number of rounds
$tbl.Rows.Item(10).Cells.Item(2).Range.Text = "anannas"             # Number of offers available
$tbl.Rows.Item(11).Cells.Item(2).Range.Text = "anannas"             # Eligibility Criteria
